$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet/tab as described in the diff (sheet name "ShearF-HW45.xpc" -> "ShearF")
$ws.Name = "ShearF"

# Append a new row of averaged-intensity data (row 16), mirroring the existing rows' layout
$row = 16

# Copy the formatting (cell style) of the A column "index" cell from the row above,
# so the new index cell (A16) keeps the same bordered/centered/bold style as A3:A15.
$ws.Cells.Item($row - 1, 1).Copy() | Out-Null
$ws.Cells.Item($row, 1).PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Cells.Item($row, 1).Value = 14
$ws.Cells.Item($row, 2).Value = "HexGrid-60degTilt5degRes"

$ws.Cells.Item($row, 3).Value = 0.9952648695185737
$ws.Cells.Item($row, 4).Value = 1.012708961054975
$ws.Cells.Item($row, 5).Value = 0.998206045989162
$ws.Cells.Item($row, 6).Value = 0.9952648695185737
$ws.Cells.Item($row, 7).Value = 1.011692470367274
$ws.Cells.Item($row, 8).Value = 0.9971603071729759
$ws.Cells.Item($row, 9).Value = 1.000689718937682
$ws.Cells.Item($row, 10).Value = 1.012708961054975
$ws.Cells.Item($row, 11).Value = 1.005457503522069
$ws.Cells.Item($row, 12).Value = 1.000361186520321
$ws.Cells.Item($row, 13).Value = 1.002620395506774
